$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for a new data row: shift the existing example rows
#        (currently rows 2-14) down by one row (to rows 3-15), without
#        using Rows.Insert() / Copy-Paste (both of which would blow the
#        existing shared formulas apart). We do this by literally moving
#        the values bottom-up.

# Row 15 is brand new, so first give it the same formatting (styles) as
# the row above it. Copy one column at a time (rather than the whole
# A14:I14 row in one go) so Excel does not also carry over / recompute a
# custom row height for the wrapped long-text cells.
foreach ($col in "A","B","C","D","E","F","G","H","I") {
    $ws.Range("${col}14").Copy()
    $ws.Range("${col}15").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ... and the constant B:F template values as the row above it.
$ws.Range("B15").Value2 = $ws.Range("B14").Value2
$ws.Range("C15").Value2 = $ws.Range("C14").Value2
$ws.Range("D15").Value2 = $ws.Range("D14").Value2
$ws.Range("E15").Value2 = $ws.Range("E14").Value2
$ws.Range("F15").Value2 = $ws.Range("F14").Value2

# Writing those multi-line template strings into the previously empty row
# makes the engine auto-grow the row height; put it back to the (shared,
# default) height used by every other row in the table.
$ws.Rows.Item(15).EntireRow.AutoFit()

# Now shift column A (the only column that actually varies row to row)
# down by one, starting from the bottom so we never overwrite a value
# before it has been copied onward.
for ($r = 14; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Range("A$dest").Value2 = $ws.Range("A$r").Value2
}

# --- 2. Populate the new first data row (row 2) with the new example.
$ws.Range("A2").Value2 = "designcourse-gsportfolio"

# --- 3. Re-assert the shared CONCAT formulas across the full, now one
#        row taller, range so the formula group covers G2:G15 / H2:H15 /
#        I2:I15 (mirrors the ref="G2:G14" -> ref="G2:G15" widening, etc.)
$ws.Range("G2:G15").Formula = "=_xlfn.CONCAT(`$B2,`$A2,`$C2,`$A2,`$D2)"
$ws.Range("H2:H15").Formula = "=_xlfn.CONCAT(`$A2,`$E2,`$A2,`$F2)"
$ws.Range("I2:I15").Formula = "=_xlfn.CONCAT(`$G2,`$H2)"

# --- 4. Update the sheet view/selection to match the post-edit state.
$ws.Range("I2").Select()
